$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rework -----------------------------------------------
# Three brand-new columns are inserted at the front (A:C) and the
# trailing "Localidade" column is replaced by a run of new columns,
# which shifts everything that used to live in A:Z and AD:AM over by
# three slots. Net effect on the surviving header cells: only text
# changes (same shared style) for most of them, while a handful of
# cells also lose their old style (no "s" attribute in the final
# sheet, i.e. the default "Normal" style).

# New leading columns (previously Concessionária / Nº da Conta / Identificador)
$ws.Range("A1").Value = "FORNECEDOR"
$ws.Range("A1").Style = "Normal"

$ws.Range("B1").Value = "CONTRATO"
$ws.Range("B1").Style = "Normal"

$ws.Range("C1").Value = "FATURA"
$ws.Range("C1").Style = "Normal"

# Old "Emissão" header becomes the unstyled "DATA EMISSÃO" column.
$ws.Range("Z1").Value = "DATA EMISSÃO"
$ws.Range("Z1").Style = "Normal"

# "Nota Fiscal" -> "NOTA FISCAL" (style unchanged).
$ws.Range("AD1").Value = "NOTA FISCAL"

# Fix a typo: "Dc_indentificador_layout" -> "Dc_identificador_layout" (style unchanged).
$ws.Range("AJ1").Value = "Dc_identificador_layout"

# Old "Localidade" header becomes the unstyled "LOCALIDADE" column.
$ws.Range("AM1").Value = "LOCALIDADE"
$ws.Range("AM1").Style = "Normal"

# --- View state ---------------------------------------------------------
# Scroll the sheet so column AD is at the left edge and select AJ2, as
# reflected in the saved sheetView.
$win = $excel.ActiveWindow
$win.ScrollColumn = 30
$win.ScrollRow = 1
$ws.Range("AJ2").Select()
